# Add notes creation for product problems and returned products during import.
#
# "Produk Retur" (sheet 1): the "Qty" column becomes a "Note" column, and the
# quantity numbers in the sample rows are replaced with free-text notes
# describing the returned product's problem.
# Also the active sheet/selection moves from "Problem Produk" (sheet2, B7)
# to "Produk Retur" (sheet1, A2).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Produk Retur"

# --- Rename the 4th column from "Qty" to "Note" and fill it with notes ----
# (changing the header cell text automatically renames the bound table
# column, e.g. Table1's tableColumn id="6")
$ws1.Range("D1").Value = "Note"
$ws1.Range("D2").Value = "Tidak Dingin"
$ws1.Range("D3").Value = "Patah"
$ws1.Range("D4").Value = "Tidak Berfungsi"

# --- Column widths on "Produk Retur": C/D no longer auto-fit ---------------
$ws1.Columns.Item(3).ColumnWidth = 12.76
$ws1.Columns.Item(4).ColumnWidth = 22.6

# --- Switch the active tab / selection --------------------------------------
# "Produk Retur" becomes the selected tab with A2 selected (previously it was
# "Problem Produk" with B7 selected).
$ws1.Activate()
[void]$ws1.Range("A2").Select()
